$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "tomorrow"
$ws.Range("B4").Value = 20
